# Applies the cryptos list update described in the commit:
# "Updated cryptos list on Tue Jan 23 22:43:06 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.466.15"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "2.224.58"
$ws.Range("E3").Value = "  -3.85%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "297.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "82.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.92%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.514"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.90%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.472"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.94%  "
$ws.Range("E10").Value = "  -3.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "29.85"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.74"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -11.15%  "
$ws.Range("D14").Value = "2.554.95"
$ws.Range("E14").Value = "  -4.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.36%  "
$ws.Range("D17").Value = "2.217.35"
$ws.Range("E17").Value = "  -4.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.718"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.14%  "
$ws.Range("D19").Value = "39.394.23"
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("E20").Value = "  -2.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "229.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.77%  "
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("E26").Value = "  -5.06%  "
$ws.Range("E27").Value = "  +1.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.49%  "
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "148.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.77%  "
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("E34").Value = "  -4.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0700"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.40%  "
$ws.Range("E36").Value = "  -4.00%  "
$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "15.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.59%  "
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.111"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0966"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.66%  "
$ws.Range("E41").Value = "  -1.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.89%  "
$ws.Range("D43").Value = "1.914.58"
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0262"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.01%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.26%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.29%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.87%  "
$ws.Range("D49").Value = "2.429.58"
$ws.Range("E49").Value = "  -4.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "88.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.30%  "
